$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1063.125
$ws.Range("I19").Value = 933.3333
$ws.Range("J19").Value = 1141
$ws.Range("K19").Value = 933.3333
$ws.Range("L19").Value = 1141
$ws.Range("M19").Value = -758.3333
$ws.Range("N19").Value = -1491

$ws.Range("H28").Value = 682.3125
$ws.Range("I28").Value = 627.8
$ws.Range("J28").Value = 1500
$ws.Range("K28").Value = 627.8
$ws.Range("L28").Value = 1500
$ws.Range("M28").Value = -142.8
$ws.Range("N28").Value = -2470

$ws.Range("H38").Value = 1082.2858
$ws.Range("I38").Value = 596.1667
$ws.Range("K38").Value = 1788.5001
$ws.Range("M38").Value = -1416.5001

$ws.Range("H127").Value = 2281.0833
$ws.Range("I127").Value = 2162.8696
$ws.Range("K127").Value = 6488.6088
$ws.Range("M127").Value = -1528.6088

$ws.Range("H132").Value = 13646.6
$ws.Range("I132").Value = 15183.5625
$ws.Range("K132").Value = 45550.6875
$ws.Range("M132").Value = -43020.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 171365.05
$ws.Range("I32").Value = 176842.8
$ws.Range("K32").Value = 176842.8
$ws.Range("M32").Value = -176555.8

$ws.Range("H62").Value = 40000
$ws.Range("J62").Value = 40000
$ws.Range("L62").Value = 40000
$ws.Range("N62").Value = -41248

$ws.Range("H65").Value = 40000
$ws.Range("J65").Value = 40000
$ws.Range("L65").Value = 120000
$ws.Range("N65").Value = -126240

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H107").Value = 65000
$ws.Range("J107").Value = 65000
$ws.Range("L107").Value = 65000
$ws.Range("N107").Value = -72680

$ws.Range("H110").Value = 965.4
$ws.Range("I110").Value = 812.2353000000001
$ws.Range("K110").Value = 812.2353000000001
$ws.Range("M110").Value = 1232.7647

$ws.Range("H122").Value = 2468.1538
$ws.Range("I122").Value = 2307.818
$ws.Range("K122").Value = 6923.454000000001
$ws.Range("M122").Value = -4473.454000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 26.25
$ws.Range("I11").Value = 26.25
$ws.Range("K11").Value = 26.25
$ws.Range("M11").Value = 113.75

$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H80").Value = 1302.5454
$ws.Range("I80").Value = 1234.2858
$ws.Range("K80").Value = 1234.2858
$ws.Range("M80").Value = -236.2858000000001

$ws.Range("H83").Value = 1302.5454
$ws.Range("I83").Value = 1234.2858
$ws.Range("K83").Value = 6171.429
$ws.Range("M83").Value = -1179.429

$ws.Range("H94").Value = 2032.7307
$ws.Range("I94").Value = 2633.1667
$ws.Range("J94").Value = 681.75
$ws.Range("K94").Value = 2633.1667
$ws.Range("L94").Value = 681.75
$ws.Range("M94").Value = -2182.1667
$ws.Range("N94").Value = -1583.75

$ws.Range("H99").Value = 5087.96
$ws.Range("I99").Value = 5518.136
$ws.Range("J99").Value = 1933.3334
$ws.Range("K99").Value = 5518.136
$ws.Range("L99").Value = 1933.3334
$ws.Range("M99").Value = -4020.136
$ws.Range("N99").Value = -4929.3334

$ws.Range("H132").Value = 125000
$ws.Range("J132").Value = 125000
$ws.Range("L132").Value = 125000
$ws.Range("N132").Value = -135120

$ws.Range("H134").Value = 4164.234
$ws.Range("I134").Value = 1889.175
$ws.Range("K134").Value = 5667.525
$ws.Range("M134").Value = -3132.525

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 27986.75
$ws.Range("J9").Value = 27986.75
$ws.Range("L9").Value = 27986.75
$ws.Range("N9").Value = -28322.75

$ws.Range("H16").Value = 91784.45
$ws.Range("I16").Value = 897
$ws.Range("J16").Value = 250837.5
$ws.Range("K16").Value = 897
$ws.Range("L16").Value = 250837.5
$ws.Range("M16").Value = -610
$ws.Range("N16").Value = -251411.5

$ws.Range("H33").Value = 243.66667
$ws.Range("I33").Value = 243.66667
$ws.Range("K33").Value = 243.66667
$ws.Range("M33").Value = 135.33333

$ws.Range("H35").Value = 3352.6
$ws.Range("I35").Value = 2256
$ws.Range("J35").Value = 4997.5
$ws.Range("K35").Value = 2256
$ws.Range("L35").Value = 4997.5
$ws.Range("M35").Value = -1962
$ws.Range("N35").Value = -5585.5

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H38").Value = 6019
$ws.Range("J38").Value = 6000
$ws.Range("L38").Value = 6000
$ws.Range("N38").Value = -6754

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H46").Value = 6019
$ws.Range("J46").Value = 6000
$ws.Range("L46").Value = 6000
$ws.Range("N46").Value = -6422

$ws.Range("H105").Value = 21430
$ws.Range("I105").Value = 25562.5
$ws.Range("J105").Value = 4900
$ws.Range("K105").Value = 25562.5
$ws.Range("L105").Value = 4900
$ws.Range("M105").Value = -23815.5
$ws.Range("N105").Value = -8394

$ws.Range("H113").Value = 91784.45
$ws.Range("I113").Value = 897
$ws.Range("J113").Value = 250837.5
$ws.Range("K113").Value = 897
$ws.Range("L113").Value = 250837.5
$ws.Range("M113").Value = 1273
$ws.Range("N113").Value = -255177.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3263
$ws.Range("I3").Value = 2924.6428
$ws.Range("J3").Value = 8000
$ws.Range("K3").Value = 8773.928400000001
$ws.Range("L3").Value = 24000
$ws.Range("M3").Value = -8661.928400000001
$ws.Range("N3").Value = -24224

$ws.Range("H11").Value = 90909270
$ws.Range("I11").Value = 168.5
$ws.Range("J11").Value = 333333540
$ws.Range("K11").Value = 505.5
$ws.Range("L11").Value = 1000000620
$ws.Range("M11").Value = -365.5
$ws.Range("N11").Value = -1000000900

$ws.Range("H99").Value = 4342.5713
$ws.Range("I99").Value = 3079.6
$ws.Range("J99").Value = 7500
$ws.Range("K99").Value = 9238.799999999999
$ws.Range("L99").Value = 22500
$ws.Range("M99").Value = -6992.799999999999
$ws.Range("N99").Value = -26992

$ws.Range("H131").Value = 2668.6316
$ws.Range("J131").Value = 2818.8572
$ws.Range("L131").Value = 8456.571599999999
$ws.Range("N131").Value = -18536.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2339.4
$ws.Range("I113").Value = 2339.4
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2339.4
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -169.4000000000001
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 2986.516
$ws.Range("J122").Value = 3329.875
$ws.Range("L122").Value = 9989.625
$ws.Range("N122").Value = -14889.625

$ws.Range("H126").Value = 2674.625
$ws.Range("I126").Value = 2474.25
$ws.Range("J126").Value = 2875
$ws.Range("K126").Value = 7422.75
$ws.Range("L126").Value = 8625
$ws.Range("M126").Value = -4952.75
$ws.Range("N126").Value = -13565

$ws.Range("H132").Value = 9535.666999999999
$ws.Range("I132").Value = 13686.6
$ws.Range("J132").Value = 3724.36
$ws.Range("K132").Value = 41059.8
$ws.Range("L132").Value = 11173.08
$ws.Range("M132").Value = -38529.8
$ws.Range("N132").Value = -16233.08

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H22").Value = 2599.5
$ws.Range("J22").Value = 2599.5
$ws.Range("L22").Value = 2599.5
$ws.Range("N22").Value = -3189.5

$ws.Range("H27").Value = 2599.5
$ws.Range("J27").Value = 2599.5
$ws.Range("L27").Value = 2599.5
$ws.Range("N27").Value = -2813.5

$ws.Range("H46").Value = 5888.222
$ws.Range("J46").Value = 5888.222
$ws.Range("L46").Value = 5888.222
$ws.Range("N46").Value = -6264.222

$ws.Range("H132").Value = 2002283
$ws.Range("I132").Value = 2502056.2
$ws.Range("J132").Value = 3189.8
$ws.Range("K132").Value = 7506168.600000001
$ws.Range("L132").Value = 9569.400000000001
$ws.Range("M132").Value = -7503638.600000001
$ws.Range("N132").Value = -14629.4

$ws.Range("H136").Value = 5955.2
$ws.Range("I136").Value = 2929.15
$ws.Range("K136").Value = 8787.450000000001
$ws.Range("M136").Value = -6237.450000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 15000
$ws.Range("I39").Value = 15000
$ws.Range("K39").Value = 15000
$ws.Range("M39").Value = -14587

$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
